$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text (the source data uses
# dotted thousand separators and fixed-format strings, not numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '66.147.96'
$ws.Range('E2').Value = '  -2.45%  '
$ws.Range('D3').Value = '3.320.78'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('D5').Value = '574.29'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').Value = '181.61'
$ws.Range('E6').Value = '  -2.58%  '
$ws.Range('D7').Value = '0.617'
$ws.Range('E7').Value = '  +3.64%  '
$ws.Range('D9').Value = '0.128'
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('D10').Value = '6.66'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('D11').Value = '0.403'
$ws.Range('E11').Value = '  -1.03%  '
$ws.Range('D12').Value = '3.894.87'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('E13').Value = '  -1.02%  '
$ws.Range('D14').Value = '26.75'
$ws.Range('E14').Value = '  -2.48%  '
$ws.Range('D15').Value = '66.222.32'
$ws.Range('E15').Value = '  -2.52%  '
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '3.308.62'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = '435.42'
$ws.Range('E18').Value = '  -2.62%  '
$ws.Range('D19').Value = '13.57'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('D20').Value = '5.66'
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('D21').Value = '7.57'
$ws.Range('E21').Value = '  -1.57%  '
$ws.Range('D22').Value = '73.37'
$ws.Range('E22').Value = '  -2.35%  '
$ws.Range('E23').Value = '  +0.35%  '
$ws.Range('D24').Value = '0.521'
$ws.Range('E24').Value = '  +1.41%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.462.50'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '0.0000116'
$ws.Range('E26').Value = '  -1.67%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '0.194'
$ws.Range('E27').Value = '  +3.29%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '9.10'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '1.95'
$ws.Range('E30').Value = '  -1.24%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '22.73'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('B32').Value = 'USDe'
$ws.Range('C32').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').Value = '5.25'
$ws.Range('E33').Value = '  -1.70%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').Value = '6.77'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = '1.22'
$ws.Range('E35').Value = '  -1.88%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.48'
$ws.Range('E36').Value = '  -2.14%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '159.66'
$ws.Range('E37').Value = '  -2.41%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = '27.62'
$ws.Range('E38').Value = '  +2.85%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '1.80'
$ws.Range('E39').Value = '  -3.72%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '2.838.15'
$ws.Range('E40').Value = '  +5.41%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').Value = '0.790'
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '4.45'
$ws.Range('E42').Value = '  -1.44%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '6.19'
$ws.Range('E43').Value = '  -2.62%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '40.52'
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = '0.0667'
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('D46').Value = '24.23'
$ws.Range('E46').Value = '  -1.53%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '2.35'
$ws.Range('E47').Value = '  -2.16%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '325.62'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0272'
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '0.102'
$ws.Range('E50').Value = '  +1.89%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = '6.15'
$ws.Range('E51').Value = '  -0.54%  '
